# Updates the LR-pairs TPM data table (Efna4-Epha2) with refreshed values
# and adds the "Resolving-Mac" sending-cluster block (rows 22-26), matching
# the new TPM computation referenced in the commit "update scripts wuth new tpm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna4"
$ws.Cells.Item(2, 3).Value = "Epha2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.974568
$ws.Cells.Item(2, 8).Value = 5.923704
$ws.Cells.Item(2, 9).Value = 0.5990695552080697
$ws.Cells.Item(2, 10).Value = 0.5990695552080698
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 20.66830833333333
$ws.Cells.Item(2, 14).Value = 62.004925
$ws.Cells.Item(2, 15).Value = 0.6755285375771634
$ws.Cells.Item(2, 16).Value = 0.6755285375771634
$ws.Cells.Item(2, 17).Value = 40.81098024913333
$ws.Cells.Item(2, 18).Value = 367.2988222422
$ws.Cells.Item(2, 19).Value = 0.404688580536709
$ws.Cells.Item(2, 20).Value = 0.4046885805367091

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna4"
$ws.Cells.Item(3, 3).Value = "Epha2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.974568
$ws.Cells.Item(3, 8).Value = 5.923704
$ws.Cells.Item(3, 9).Value = 0.5990695552080697
$ws.Cells.Item(3, 10).Value = 0.5990695552080698
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.156330666666667
$ws.Cells.Item(3, 14).Value = 6.468992
$ws.Cells.Item(3, 15).Value = 0.07047809033489469
$ws.Cells.Item(3, 16).Value = 0.07047809033489467
$ws.Cells.Item(3, 17).Value = 4.257821531818667
$ws.Cells.Item(3, 18).Value = 38.320393786368
$ws.Cells.Item(3, 19).Value = 0.04222127822883952
$ws.Cells.Item(3, 20).Value = 0.04222127822883952

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna4"
$ws.Cells.Item(4, 3).Value = "Epha2"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.974568
$ws.Cells.Item(4, 8).Value = 5.923704
$ws.Cells.Item(4, 9).Value = 0.5990695552080697
$ws.Cells.Item(4, 10).Value = 0.5990695552080698
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.683564
$ws.Cells.Item(4, 14).Value = 5.050692
$ws.Cells.Item(4, 15).Value = 0.05502605769642779
$ws.Cells.Item(4, 16).Value = 0.05502605769642779
$ws.Cells.Item(4, 17).Value = 3.324311600351999
$ws.Cells.Item(4, 18).Value = 29.918804403168
$ws.Cells.Item(4, 19).Value = 0.03296443590905258
$ws.Cells.Item(4, 20).Value = 0.03296443590905258

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Efna4"
$ws.Cells.Item(5, 3).Value = "Epha2"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.974568
$ws.Cells.Item(5, 8).Value = 5.923704
$ws.Cells.Item(5, 9).Value = 0.5990695552080697
$ws.Cells.Item(5, 10).Value = 0.5990695552080698
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 5.278649666666666
$ws.Cells.Item(5, 14).Value = 15.835949
$ws.Cells.Item(5, 15).Value = 0.1725288026574751
$ws.Cells.Item(5, 16).Value = 0.1725288026574751
$ws.Cells.Item(5, 17).Value = 10.42305271501067
$ws.Cells.Item(5, 18).Value = 93.807474435096
$ws.Cells.Item(5, 19).Value = 0.1033567530685944
$ws.Cells.Item(5, 20).Value = 0.1033567530685944

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Efna4"
$ws.Cells.Item(6, 3).Value = "Epha2"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.974568
$ws.Cells.Item(6, 8).Value = 5.923704
$ws.Cells.Item(6, 9).Value = 0.5990695552080697
$ws.Cells.Item(6, 10).Value = 0.5990695552080698
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.8089063333333334
$ws.Cells.Item(6, 14).Value = 2.426719
$ws.Cells.Item(6, 15).Value = 0.02643851173403914
$ws.Cells.Item(6, 16).Value = 0.02643851173403913
$ws.Cells.Item(6, 17).Value = 1.597240560797333
$ws.Cells.Item(6, 18).Value = 14.375165047176
$ws.Cells.Item(6, 19).Value = 0.01583850746487416
$ws.Cells.Item(6, 20).Value = 0.01583850746487416

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna4"
$ws.Cells.Item(7, 3).Value = "Epha2"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.020259
$ws.Cells.Item(7, 8).Value = 3.060777
$ws.Cells.Item(7, 9).Value = 0.309539152527049
$ws.Cells.Item(7, 10).Value = 0.309539152527049
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 20.66830833333333
$ws.Cells.Item(7, 14).Value = 62.004925
$ws.Cells.Item(7, 15).Value = 0.6755285375771634
$ws.Cells.Item(7, 16).Value = 0.6755285375771634
$ws.Cells.Item(7, 17).Value = 21.08702759185833
$ws.Cells.Item(7, 18).Value = 189.783248326725
$ws.Cells.Item(7, 19).Value = 0.2091025310294719
$ws.Cells.Item(7, 20).Value = 0.2091025310294719

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Efna4"
$ws.Cells.Item(8, 3).Value = "Epha2"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.020259
$ws.Cells.Item(8, 8).Value = 3.060777
$ws.Cells.Item(8, 9).Value = 0.309539152527049
$ws.Cells.Item(8, 10).Value = 0.309539152527049
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.156330666666667
$ws.Cells.Item(8, 14).Value = 6.468992
$ws.Cells.Item(8, 15).Value = 0.07047809033489469
$ws.Cells.Item(8, 16).Value = 0.07047809033489467
$ws.Cells.Item(8, 17).Value = 2.200015769642667
$ws.Cells.Item(8, 18).Value = 19.800141926784
$ws.Cells.Item(8, 19).Value = 0.02181572835398811
$ws.Cells.Item(8, 20).Value = 0.0218157283539881

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Efna4"
$ws.Cells.Item(9, 3).Value = "Epha2"
$ws.Cells.Item(9, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.020259
$ws.Cells.Item(9, 8).Value = 3.060777
$ws.Cells.Item(9, 9).Value = 0.309539152527049
$ws.Cells.Item(9, 10).Value = 0.309539152527049
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.683564
$ws.Cells.Item(9, 14).Value = 5.050692
$ws.Cells.Item(9, 15).Value = 0.05502605769642779
$ws.Cells.Item(9, 16).Value = 0.05502605769642779
$ws.Cells.Item(9, 17).Value = 1.717671323076
$ws.Cells.Item(9, 18).Value = 15.459041907684
$ws.Cells.Item(9, 19).Value = 0.01703271926625676
$ws.Cells.Item(9, 20).Value = 0.01703271926625676

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Efna4"
$ws.Cells.Item(10, 3).Value = "Epha2"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.020259
$ws.Cells.Item(10, 8).Value = 3.060777
$ws.Cells.Item(10, 9).Value = 0.309539152527049
$ws.Cells.Item(10, 10).Value = 0.309539152527049
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 5.278649666666666
$ws.Cells.Item(10, 14).Value = 15.835949
$ws.Cells.Item(10, 15).Value = 0.1725288026574751
$ws.Cells.Item(10, 16).Value = 0.1725288026574751
$ws.Cells.Item(10, 17).Value = 5.385589830263666
$ws.Cells.Item(10, 18).Value = 48.470308472373
$ws.Cells.Item(10, 19).Value = 0.05340441936110132
$ws.Cells.Item(10, 20).Value = 0.05340441936110131

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Efna4"
$ws.Cells.Item(11, 3).Value = "Epha2"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.020259
$ws.Cells.Item(11, 8).Value = 3.060777
$ws.Cells.Item(11, 9).Value = 0.309539152527049
$ws.Cells.Item(11, 10).Value = 0.309539152527049
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.8089063333333334
$ws.Cells.Item(11, 14).Value = 2.426719
$ws.Cells.Item(11, 15).Value = 0.02643851173403914
$ws.Cells.Item(11, 16).Value = 0.02643851173403913
$ws.Cells.Item(11, 17).Value = 0.8252939667403334
$ws.Cells.Item(11, 18).Value = 7.427645700663001
$ws.Cells.Item(11, 19).Value = 0.008183754516230915
$ws.Cells.Item(11, 20).Value = 0.008183754516230914

# Row 12
$ws.Cells.Item(12, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 2).Value = "Efna4"
$ws.Cells.Item(12, 3).Value = "Epha2"
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.077601
$ws.Cells.Item(12, 8).Value = 0.232803
$ws.Cells.Item(12, 9).Value = 0.02354357842003994
$ws.Cells.Item(12, 10).Value = 0.02354357842003994
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 20.66830833333333
$ws.Cells.Item(12, 14).Value = 62.004925
$ws.Cells.Item(12, 15).Value = 0.6755285375771634
$ws.Cells.Item(12, 16).Value = 0.6755285375771634
$ws.Cells.Item(12, 17).Value = 1.603881394975
$ws.Cells.Item(12, 18).Value = 14.434932554775
$ws.Cells.Item(12, 19).Value = 0.01590435909942284
$ws.Cells.Item(12, 20).Value = 0.01590435909942284

# Row 13
$ws.Cells.Item(13, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13, 2).Value = "Efna4"
$ws.Cells.Item(13, 3).Value = "Epha2"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.077601
$ws.Cells.Item(13, 8).Value = 0.232803
$ws.Cells.Item(13, 9).Value = 0.02354357842003994
$ws.Cells.Item(13, 10).Value = 0.02354357842003994
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.156330666666667
$ws.Cells.Item(13, 14).Value = 6.468992
$ws.Cells.Item(13, 15).Value = 0.07047809033489469
$ws.Cells.Item(13, 16).Value = 0.07047809033489467
$ws.Cells.Item(13, 17).Value = 0.167333416064
$ws.Cells.Item(13, 18).Value = 1.506000744576
$ws.Cells.Item(13, 19).Value = 0.001659306446694252
$ws.Cells.Item(13, 20).Value = 0.001659306446694252

# Row 14
$ws.Cells.Item(14, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(14, 2).Value = "Efna4"
$ws.Cells.Item(14, 3).Value = "Epha2"
$ws.Cells.Item(14, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.077601
$ws.Cells.Item(14, 8).Value = 0.232803
$ws.Cells.Item(14, 9).Value = 0.02354357842003994
$ws.Cells.Item(14, 10).Value = 0.02354357842003994
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 1.683564
$ws.Cells.Item(14, 14).Value = 5.050692
$ws.Cells.Item(14, 15).Value = 0.05502605769642779
$ws.Cells.Item(14, 16).Value = 0.05502605769642779
$ws.Cells.Item(14, 17).Value = 0.130646249964
$ws.Cells.Item(14, 18).Value = 1.175816249676
$ws.Cells.Item(14, 19).Value = 0.00129551030452149
$ws.Cells.Item(14, 20).Value = 0.00129551030452149

# Row 15
$ws.Cells.Item(15, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(15, 2).Value = "Efna4"
$ws.Cells.Item(15, 3).Value = "Epha2"
$ws.Cells.Item(15, 4).Value = "MuSCs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.077601
$ws.Cells.Item(15, 8).Value = 0.232803
$ws.Cells.Item(15, 9).Value = 0.02354357842003994
$ws.Cells.Item(15, 10).Value = 0.02354357842003994
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 5.278649666666666
$ws.Cells.Item(15, 14).Value = 15.835949
$ws.Cells.Item(15, 15).Value = 0.1725288026574751
$ws.Cells.Item(15, 16).Value = 0.1725288026574751
$ws.Cells.Item(15, 17).Value = 0.409628492783
$ws.Cells.Item(15, 18).Value = 3.686656435047
$ws.Cells.Item(15, 19).Value = 0.00406194539508186
$ws.Cells.Item(15, 20).Value = 0.00406194539508186

# Row 16
$ws.Cells.Item(16, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 2).Value = "Efna4"
$ws.Cells.Item(16, 3).Value = "Epha2"
$ws.Cells.Item(16, 4).Value = "Resolving-Mac"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.077601
$ws.Cells.Item(16, 8).Value = 0.232803
$ws.Cells.Item(16, 9).Value = 0.02354357842003994
$ws.Cells.Item(16, 10).Value = 0.02354357842003994
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.8089063333333334
$ws.Cells.Item(16, 14).Value = 2.426719
$ws.Cells.Item(16, 15).Value = 0.02643851173403914
$ws.Cells.Item(16, 16).Value = 0.02643851173403913
$ws.Cells.Item(16, 17).Value = 0.062771940373
$ws.Cells.Item(16, 18).Value = 0.5649474633570001
$ws.Cells.Item(16, 19).Value = 0.0006224571743194965
$ws.Cells.Item(16, 20).Value = 0.0006224571743194965

# Row 17
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Efna4"
$ws.Cells.Item(17, 3).Value = "Epha2"
$ws.Cells.Item(17, 4).Value = "ECs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.1837383333333333
$ws.Cells.Item(17, 8).Value = 0.551215
$ws.Cells.Item(17, 9).Value = 0.05574487261247628
$ws.Cells.Item(17, 10).Value = 0.05574487261247628
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 20.66830833333333
$ws.Cells.Item(17, 14).Value = 62.004925
$ws.Cells.Item(17, 15).Value = 0.6755285375771634
$ws.Cells.Item(17, 16).Value = 0.6755285375771634
$ws.Cells.Item(17, 17).Value = 3.797560525986111
$ws.Cells.Item(17, 18).Value = 34.178044733875
$ws.Cells.Item(17, 19).Value = 0.03765725227333137
$ws.Cells.Item(17, 20).Value = 0.03765725227333137

# Row 18
$ws.Cells.Item(18, 1).Value = "MuSCs"
$ws.Cells.Item(18, 2).Value = "Efna4"
$ws.Cells.Item(18, 3).Value = "Epha2"
$ws.Cells.Item(18, 4).Value = "FAPs"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 0.1837383333333333
$ws.Cells.Item(18, 8).Value = 0.551215
$ws.Cells.Item(18, 9).Value = 0.05574487261247628
$ws.Cells.Item(18, 10).Value = 0.05574487261247628
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 2.156330666666667
$ws.Cells.Item(18, 14).Value = 6.468992
$ws.Cells.Item(18, 15).Value = 0.07047809033489469
$ws.Cells.Item(18, 16).Value = 0.07047809033489467
$ws.Cells.Item(18, 17).Value = 0.3962006028088889
$ws.Cells.Item(18, 18).Value = 3.56580542528
$ws.Cells.Item(18, 19).Value = 0.0039287921676893
$ws.Cells.Item(18, 20).Value = 0.003928792167689299

# Row 19
$ws.Cells.Item(19, 1).Value = "MuSCs"
$ws.Cells.Item(19, 2).Value = "Efna4"
$ws.Cells.Item(19, 3).Value = "Epha2"
$ws.Cells.Item(19, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 0.1837383333333333
$ws.Cells.Item(19, 8).Value = 0.551215
$ws.Cells.Item(19, 9).Value = 0.05574487261247628
$ws.Cells.Item(19, 10).Value = 0.05574487261247628
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 1.683564
$ws.Cells.Item(19, 14).Value = 5.050692
$ws.Cells.Item(19, 15).Value = 0.05502605769642779
$ws.Cells.Item(19, 16).Value = 0.05502605769642779
$ws.Cells.Item(19, 17).Value = 0.30933524342
$ws.Cells.Item(19, 18).Value = 2.78401719078
$ws.Cells.Item(19, 19).Value = 0.003067420576654137
$ws.Cells.Item(19, 20).Value = 0.003067420576654137

# Row 20
$ws.Cells.Item(20, 1).Value = "MuSCs"
$ws.Cells.Item(20, 2).Value = "Efna4"
$ws.Cells.Item(20, 3).Value = "Epha2"
$ws.Cells.Item(20, 4).Value = "MuSCs"
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 0.1837383333333333
$ws.Cells.Item(20, 8).Value = 0.551215
$ws.Cells.Item(20, 9).Value = 0.05574487261247628
$ws.Cells.Item(20, 10).Value = 0.05574487261247628
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 5.278649666666666
$ws.Cells.Item(20, 14).Value = 15.835949
$ws.Cells.Item(20, 15).Value = 0.1725288026574751
$ws.Cells.Item(20, 16).Value = 0.1725288026574751
$ws.Cells.Item(20, 17).Value = 0.9698902920038889
$ws.Cells.Item(20, 18).Value = 8.729012628034999
$ws.Cells.Item(20, 19).Value = 0.009617596126124008
$ws.Cells.Item(20, 20).Value = 0.009617596126124006

# Row 21
$ws.Cells.Item(21, 1).Value = "MuSCs"
$ws.Cells.Item(21, 2).Value = "Efna4"
$ws.Cells.Item(21, 3).Value = "Epha2"
$ws.Cells.Item(21, 4).Value = "Resolving-Mac"
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 0.1837383333333333
$ws.Cells.Item(21, 8).Value = 0.551215
$ws.Cells.Item(21, 9).Value = 0.05574487261247628
$ws.Cells.Item(21, 10).Value = 0.05574487261247628
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 0.8089063333333334
$ws.Cells.Item(21, 14).Value = 2.426719
$ws.Cells.Item(21, 15).Value = 0.02643851173403914
$ws.Cells.Item(21, 16).Value = 0.02643851173403913
$ws.Cells.Item(21, 17).Value = 0.1486271015094445
$ws.Cells.Item(21, 18).Value = 1.337643913585
$ws.Cells.Item(21, 19).Value = 0.001473811468677471
$ws.Cells.Item(21, 20).Value = 0.001473811468677471

# Row 22
$ws.Cells.Item(22, 1).Value = "Resolving-Mac"
$ws.Cells.Item(22, 2).Value = "Efna4"
$ws.Cells.Item(22, 3).Value = "Epha2"
$ws.Cells.Item(22, 4).Value = "ECs"
$ws.Cells.Item(22, 5).Value = 1
$ws.Cells.Item(22, 6).Value = 0.3333333333333333
$ws.Cells.Item(22, 7).Value = 0.03989166666666667
$ws.Cells.Item(22, 8).Value = 0.119675
$ws.Cells.Item(22, 9).Value = 0.01210284123236505
$ws.Cells.Item(22, 10).Value = 0.01210284123236505
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 12).Value = 1
$ws.Cells.Item(22, 13).Value = 20.66830833333333
$ws.Cells.Item(22, 14).Value = 62.004925
$ws.Cells.Item(22, 15).Value = 0.6755285375771634
$ws.Cells.Item(22, 16).Value = 0.6755285375771634
$ws.Cells.Item(22, 17).Value = 0.8244932665972221
$ws.Cells.Item(22, 18).Value = 7.420439399375001
$ws.Cells.Item(22, 19).Value = 0.008175814638228153
$ws.Cells.Item(22, 20).Value = 0.008175814638228154

# Row 23
$ws.Cells.Item(23, 1).Value = "Resolving-Mac"
$ws.Cells.Item(23, 2).Value = "Efna4"
$ws.Cells.Item(23, 3).Value = "Epha2"
$ws.Cells.Item(23, 4).Value = "FAPs"
$ws.Cells.Item(23, 5).Value = 1
$ws.Cells.Item(23, 6).Value = 0.3333333333333333
$ws.Cells.Item(23, 7).Value = 0.03989166666666667
$ws.Cells.Item(23, 8).Value = 0.119675
$ws.Cells.Item(23, 9).Value = 0.01210284123236505
$ws.Cells.Item(23, 10).Value = 0.01210284123236505
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 12).Value = 1
$ws.Cells.Item(23, 13).Value = 2.156330666666667
$ws.Cells.Item(23, 14).Value = 6.468992
$ws.Cells.Item(23, 15).Value = 0.07047809033489469
$ws.Cells.Item(23, 16).Value = 0.07047809033489467
$ws.Cells.Item(23, 17).Value = 0.08601962417777778
$ws.Cells.Item(23, 18).Value = 0.7741766176
$ws.Cells.Item(23, 19).Value = 0.0008529851376835118
$ws.Cells.Item(23, 20).Value = 0.0008529851376835117

# Row 24
$ws.Cells.Item(24, 1).Value = "Resolving-Mac"
$ws.Cells.Item(24, 2).Value = "Efna4"
$ws.Cells.Item(24, 3).Value = "Epha2"
$ws.Cells.Item(24, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(24, 5).Value = 1
$ws.Cells.Item(24, 6).Value = 0.3333333333333333
$ws.Cells.Item(24, 7).Value = 0.03989166666666667
$ws.Cells.Item(24, 8).Value = 0.119675
$ws.Cells.Item(24, 9).Value = 0.01210284123236505
$ws.Cells.Item(24, 10).Value = 0.01210284123236505
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 12).Value = 1
$ws.Cells.Item(24, 13).Value = 1.683564
$ws.Cells.Item(24, 14).Value = 5.050692
$ws.Cells.Item(24, 15).Value = 0.05502605769642779
$ws.Cells.Item(24, 16).Value = 0.05502605769642779
$ws.Cells.Item(24, 17).Value = 0.0671601739
$ws.Cells.Item(24, 18).Value = 0.6044415651
$ws.Cells.Item(24, 19).Value = 0.0006659716399428242
$ws.Cells.Item(24, 20).Value = 0.0006659716399428243

# Row 25
$ws.Cells.Item(25, 1).Value = "Resolving-Mac"
$ws.Cells.Item(25, 2).Value = "Efna4"
$ws.Cells.Item(25, 3).Value = "Epha2"
$ws.Cells.Item(25, 4).Value = "MuSCs"
$ws.Cells.Item(25, 5).Value = 1
$ws.Cells.Item(25, 6).Value = 0.3333333333333333
$ws.Cells.Item(25, 7).Value = 0.03989166666666667
$ws.Cells.Item(25, 8).Value = 0.119675
$ws.Cells.Item(25, 9).Value = 0.01210284123236505
$ws.Cells.Item(25, 10).Value = 0.01210284123236505
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(25, 13).Value = 5.278649666666666
$ws.Cells.Item(25, 14).Value = 15.835949
$ws.Cells.Item(25, 15).Value = 0.1725288026574751
$ws.Cells.Item(25, 16).Value = 0.1725288026574751
$ws.Cells.Item(25, 17).Value = 0.2105741329527778
$ws.Cells.Item(25, 18).Value = 1.895167196575
$ws.Cells.Item(25, 19).Value = 0.002088088706573461
$ws.Cells.Item(25, 20).Value = 0.002088088706573461

# Row 26
$ws.Cells.Item(26, 1).Value = "Resolving-Mac"
$ws.Cells.Item(26, 2).Value = "Efna4"
$ws.Cells.Item(26, 3).Value = "Epha2"
$ws.Cells.Item(26, 4).Value = "Resolving-Mac"
$ws.Cells.Item(26, 5).Value = 1
$ws.Cells.Item(26, 6).Value = 0.3333333333333333
$ws.Cells.Item(26, 7).Value = 0.03989166666666667
$ws.Cells.Item(26, 8).Value = 0.119675
$ws.Cells.Item(26, 9).Value = 0.01210284123236505
$ws.Cells.Item(26, 10).Value = 0.01210284123236505
$ws.Cells.Item(26, 11).Value = 3
$ws.Cells.Item(26, 12).Value = 1
$ws.Cells.Item(26, 13).Value = 0.8089063333333334
$ws.Cells.Item(26, 14).Value = 2.426719
$ws.Cells.Item(26, 15).Value = 0.02643851173403914
$ws.Cells.Item(26, 16).Value = 0.02643851173403913
$ws.Cells.Item(26, 17).Value = 0.03226862181388889
$ws.Cells.Item(26, 18).Value = 0.290417596325
$ws.Cells.Item(26, 19).Value = 0.0003199811099370959
$ws.Cells.Item(26, 20).Value = 0.0003199811099370959
